$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startYear = 1451
$startPop = 93401
$startRow = 53
$count = 20

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $year = $startYear + $i
    $pop = $startPop - ($i * 11)

    $cellA = $ws.Cells.Item($row, 1)
    $cellB = $ws.Cells.Item($row, 2)

    # Carry the same row-style used by the last existing data row (52)
    # down into the newly appended rows.
    $ws.Cells.Item(52, 1).Copy()
    $cellA.PasteSpecial(-4122)
    $ws.Cells.Item(52, 2).Copy()
    $cellB.PasteSpecial(-4122)

    $cellA.Value = $year
    $cellB.Value = $pop

    $ws.Rows.Item($row).RowHeight = 12.8
}

$excel.CutCopyMode = 0

$ws.Application.ActiveWindow.ScrollRow = 35
$ws.Range("D68").Select() | Out-Null
